# Apply cryptos.xlsx data refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.238.65'
$ws.Cells.Item(2, 5).Value = '  -0.02%  '
$ws.Cells.Item(3, 4).Value = '1.869.95'
$ws.Cells.Item(3, 5).Value = '  +0.22%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$ws.Cells.Item(5, 4).Value = '''0.7105'
$ws.Cells.Item(5, 5).Value = '  -1.73%  '
$ws.Cells.Item(6, 4).Value = '''241.33'
$ws.Cells.Item(6, 5).Value = '  +0.21%  '
$ws.Cells.Item(7, 4).Value = '''1.001'
$ws.Cells.Item(7, 5).Value = '  +0.00%  '
$ws.Cells.Item(8, 4).Value = '''0.3104'
$ws.Cells.Item(8, 5).Value = '  +0.54%  '
$ws.Cells.Item(9, 4).Value = '''0.07684'
$ws.Cells.Item(9, 5).Value = '  -1.96%  '
$ws.Cells.Item(10, 4).Value = '''24.96'
$ws.Cells.Item(10, 5).Value = '  -1.26%  '
$ws.Cells.Item(11, 4).Value = '''0.08352'
$ws.Cells.Item(11, 5).Value = '  +1.24%  '
$ws.Cells.Item(12, 4).Value = '1.885.26'
$ws.Cells.Item(12, 5).Value = '  +1.08%  '
$ws.Cells.Item(13, 4).Value = '''5.208'
$ws.Cells.Item(13, 5).Value = '  -0.54%  '
$ws.Cells.Item(14, 4).Value = '''0.7093'
$ws.Cells.Item(14, 5).Value = '  -1.75%  '
$ws.Cells.Item(15, 4).Value = '''91.13'
$ws.Cells.Item(16, 4).Value = '29.265.51'
$ws.Cells.Item(16, 5).Value = '  -0.08%  '
$ws.Cells.Item(17, 4).Value = '''0.000008247'
$ws.Cells.Item(17, 5).Value = '  +5.62%  '
$ws.Cells.Item(18, 4).Value = '''5.925'
$ws.Cells.Item(18, 5).Value = '  +1.13%  '
$ws.Cells.Item(19, 4).Value = '''242.06'
$ws.Cells.Item(19, 5).Value = '  -0.71%  '
$ws.Cells.Item(20, 4).Value = '2.131.39'
$ws.Cells.Item(20, 5).Value = '  +0.42%  '
$ws.Cells.Item(21, 4).Value = '''13.15'
$ws.Cells.Item(22, 4).Value = '''0.9995'
$ws.Cells.Item(22, 5).Value = '  -0.08%  '
$ws.Cells.Item(23, 4).Value = '''7.823'
$ws.Cells.Item(23, 5).Value = '  -1.86%  '
$ws.Cells.Item(24, 5).Value = '  -0.05%  '
$ws.Cells.Item(25, 4).Value = '''0.1624'
$ws.Cells.Item(25, 5).Value = '  +1.81%  '
$ws.Cells.Item(26, 4).Value = '''163.24'
$ws.Cells.Item(26, 5).Value = '  +0.85%  '
$ws.Cells.Item(27, 4).Value = '''8.995'
$ws.Cells.Item(27, 5).Value = '  +0.40%  '
$ws.Cells.Item(28, 4).Value = '''18.48'
$ws.Cells.Item(28, 5).Value = '  +1.39%  '
$ws.Cells.Item(29, 4).Value = '''1.504'
$ws.Cells.Item(29, 5).Value = '  +0.77%  '
$ws.Cells.Item(30, 5).Value = '  +0.15%  '
$ws.Cells.Item(31, 4).Value = '''4.332'
$ws.Cells.Item(31, 5).Value = '  +5.56%  '
$ws.Cells.Item(32, 4).Value = '''1.278'
$ws.Cells.Item(32, 5).Value = '  -4.85%  '
$ws.Cells.Item(33, 4).Value = '''0.05234'
$ws.Cells.Item(33, 5).Value = '  +0.77%  '
$ws.Cells.Item(34, 5).Value = '  -0.49%  '
$ws.Cells.Item(35, 4).Value = '''0.7500'
$ws.Cells.Item(35, 5).Value = '  +2.97%  '
$ws.Cells.Item(36, 4).Value = '''1.170'
$ws.Cells.Item(36, 5).Value = '  -1.24%  '
$ws.Cells.Item(37, 4).Value = '''2.683'
$ws.Cells.Item(37, 5).Value = '  +0.19%  '
$ws.Cells.Item(38, 4).Value = '''0.01854'
$ws.Cells.Item(38, 5).Value = '  -0.05%  '
$ws.Cells.Item(39, 4).Value = '''2.715'
$ws.Cells.Item(39, 5).Value = '  +0.58%  '
$ws.Cells.Item(40, 4).Value = '1.151.65'
$ws.Cells.Item(40, 5).Value = '  -1.97%  '
$ws.Cells.Item(41, 4).Value = '''6.357'
$ws.Cells.Item(41, 5).Value = '  +4.02%  '
$ws.Cells.Item(42, 4).Value = '''72.97'
$ws.Cells.Item(42, 5).Value = '  +0.76%  '
$ws.Cells.Item(43, 4).Value = '''0.8855'
$ws.Cells.Item(43, 5).Value = '  -2.00%  '
$ws.Cells.Item(44, 4).Value = '''104.70'
$ws.Cells.Item(44, 5).Value = '  +2.80%  '
$ws.Cells.Item(45, 4).Value = '''1.000'
$ws.Cells.Item(45, 5).Value = '  +0.01%  '
$ws.Cells.Item(46, 4).Value = '2.027.08'
$ws.Cells.Item(46, 5).Value = '  +0.90%  '
$ws.Cells.Item(47, 4).Value = '''0.5183'
$ws.Cells.Item(47, 5).Value = '  -1.88%  '
$ws.Cells.Item(48, 4).Value = '''1.791'
$ws.Cells.Item(48, 5).Value = '  +0.65%  '
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(49, 4).Value = '''9.353'
$ws.Cells.Item(49, 5).Value = '  +0.60%  '
$ws.Cells.Item(50, 2).Value = 'TheSandbox'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(50, 4).Value = '''0.4284'
$ws.Cells.Item(50, 5).Value = '  +0.45%  '
$ws.Cells.Item(51, 2).Value = 'Frax'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(51, 4).Value = '''0.9985'
$ws.Cells.Item(51, 5).Value = '  +0.14%  '
